# Generate Report for Handback
# This script marks the two localization jobs (zh-cn, de-de) as handed back:
#  - Overview sheet: status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both locale columns.
#  - zh-cn / de-de sheets: status column updated the same way, and the
#    "Latest Target File" / "Latest Handback File" / "Latest Handback
#    DateTime" columns (I, J, K) are now populated for both data rows.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdFile1 = "1df68d2c-6b5c-4aad-92ca-a6905a981c51.md"
$mdFile2 = "4e7f92f3-e627-4b01-898d-6ddf5a59070a.md"
$mdUrl1  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58f85dbf7513e29d81137e73c29e2ee5826cd808/e2e/1df68d2c-6b5c-4aad-92ca-a6905a981c51.md"
$mdUrl2  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58f85dbf7513e29d81137e73c29e2ee5826cd808/e2e/4e7f92f3-e627-4b01-898d-6ddf5a59070a.md"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777050018311
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777050018311

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsZh.Range("I2").Value = $mdFile1
$wsZh.Range("J2").Value = "1df68d2c-6b5c-4aad-92ca-a6905a981c51.d4cf9d6eef8b43e3992d4f1c93b6fa6487e9bfc0.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-10-14 08:29:58"

$wsZh.Range("I3").Value = $mdFile2
$wsZh.Range("J3").Value = "4e7f92f3-e627-4b01-898d-6ddf5a59070a.dc18530527bf222931cef2a3c835ea8d5e2dee90.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-10-14 08:29:58"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdFile2)

$wsZh.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$wsDe.Range("I2").Value = $mdFile1
$wsDe.Range("J2").Value = "1df68d2c-6b5c-4aad-92ca-a6905a981c51.d4cf9d6eef8b43e3992d4f1c93b6fa6487e9bfc0.de-de.xlf"
$wsDe.Range("K2").Value = "2016-10-14 08:30:17"

$wsDe.Range("I3").Value = $mdFile2
$wsDe.Range("J3").Value = "4e7f92f3-e627-4b01-898d-6ddf5a59070a.dc18530527bf222931cef2a3c835ea8d5e2dee90.de-de.xlf"
$wsDe.Range("K3").Value = "2016-10-14 08:30:17"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdFile2)

$wsDe.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

Write-Host "Handback report generated"
